$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 204, shifting existing rows 204:312 down to 205:313.
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new data record
# (constant columns copied from the surrounding rows; changed columns
# per the new record).
$ws.Cells.Item(204, 1).Value = 11
$ws.Cells.Item(204, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value = "Bíobío"
$ws.Cells.Item(204, 4).Value = 44917
$ws.Cells.Item(204, 5).Value = 8
$ws.Cells.Item(204, 6).Value = 100114013
$ws.Cells.Item(204, 7).Value = "Zanahoria"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 1000
$ws.Cells.Item(204, 11).Value = 6500
$ws.Cells.Item(204, 12).Value = 7000
$ws.Cells.Item(204, 13).Value = 6750
$ws.Cells.Item(204, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(204, 15).Value = "Región de Ñuble"
$ws.Cells.Item(204, 16).Value = 338
$ws.Cells.Item(204, 17).Value = 20
$ws.Cells.Item(204, 18).Value = "Hortaliza"
